# Romanian COM mapping scheme: simplify residential/commercial building
# classes per local expert feedback. The "ER+ETC" (Earthen/Rubble masonry
# + Earthquake-resistant timber construction) taxonomy label is replaced
# with "MUR+ADO" (Unreinforced masonry + Adobe) across all three vulnerability
# lists (Offices, Trade, Hotels), and the sheet is reformatted for
# readability (wider columns, wrapped/tall row for the long strings).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Content fix: ER+ETC -> MUR+ADO (one occurrence per column) ---
$ws.Range("B2").Value = $ws.Range("B2").Text.Replace("ER+ETC", "MUR+ADO")
$ws.Range("C2").Value = $ws.Range("C2").Text.Replace("ER+ETC", "MUR+ADO")
$ws.Range("D2").Value = $ws.Range("D2").Text.Replace("ER+ETC", "MUR+ADO")

# --- 2. Widen the columns so the multi-line labels are readable ---
$ws.Columns.Item(2).ColumnWidth = 39.166666666666664   # -> stored width 40
$ws.Columns.Item(3).ColumnWidth = 44.166666666666664   # -> stored width 45
$ws.Columns.Item(4).ColumnWidth = 45.166666666666664   # -> stored width 46

# --- 3. Let row 2 grow tall and wrap the long text blocks ---
$ws.Range("B2:D2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 409

# --- 4. Restore the view: scroll to row 2, select cell F2 ---
$ws.Range("F2").Select() | Out-Null

Write-Output "Applied Romania COM mapping scheme update"
